$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new line names to the table (line7, line8) in rows 16 and 17
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "line7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "line8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true

# Update existing values for rows 8-15 (C, D, E columns)
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 3).Value = 16

$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 3).Value = 10

$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false
